# Update "想去人数" (F column) figures across the workbook's sheets to
# reflect the regenerated data (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (row, new value) updates for column F.
$updates = @{
    "展览" = @(
        @{ Row = 2;  Value = 2676 },
        @{ Row = 3;  Value = 1029 },
        @{ Row = 4;  Value = 19216 },
        @{ Row = 6;  Value = 2151 },
        @{ Row = 7;  Value = 726 },
        @{ Row = 11; Value = 221 },
        @{ Row = 14; Value = 346 },
        @{ Row = 16; Value = 243 },
        @{ Row = 18; Value = 165 }
    )
    "演出" = @(
        @{ Row = 3; Value = 32 },
        @{ Row = 4; Value = 6 },
        @{ Row = 5; Value = 11 },
        @{ Row = 7; Value = 261 },
        @{ Row = 8; Value = 120 }
    )
    "本地生活" = @(
        @{ Row = 2; Value = 5960 },
        @{ Row = 3; Value = 619 }
    )
    "全部类型" = @(
        @{ Row = 2;  Value = 5960 },
        @{ Row = 3;  Value = 619 },
        @{ Row = 6;  Value = 32 },
        @{ Row = 7;  Value = 2676 },
        @{ Row = 8;  Value = 1029 },
        @{ Row = 9;  Value = 19216 },
        @{ Row = 10; Value = 6 },
        @{ Row = 11; Value = 11 },
        @{ Row = 14; Value = 261 },
        @{ Row = 15; Value = 2151 },
        @{ Row = 16; Value = 726 },
        @{ Row = 17; Value = 120 },
        @{ Row = 21; Value = 221 },
        @{ Row = 27; Value = 346 },
        @{ Row = 30; Value = 243 },
        @{ Row = 34; Value = 165 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Cells.Item($entry.Row, 6).Value = $entry.Value
    }
}
